# Update "想去人数" (number of people interested) figures in the
# "展览" (Exhibition) sheet and the corresponding rows in the
# "全部类型" (All Types) sheet, per the latest data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibition listing) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 54
$wsExhibit.Range("F8").Value = 80
$wsExhibit.Range("F9").Value = 8518
$wsExhibit.Range("F11").Value = 315
$wsExhibit.Range("F12").Value = 1134
$wsExhibit.Range("F13").Value = 924
$wsExhibit.Range("F14").Value = 84
$wsExhibit.Range("F16").Value = 224
$wsExhibit.Range("F17").Value = 199
$wsExhibit.Range("F18").Value = 60
$wsExhibit.Range("F19").Value = 227
$wsExhibit.Range("F20").Value = 965

# --- Sheet "全部类型" (all types, combined listing) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 54
$wsAll.Range("F10").Value = 80
$wsAll.Range("F11").Value = 8518
$wsAll.Range("F13").Value = 315
$wsAll.Range("F14").Value = 1134
$wsAll.Range("F15").Value = 924
$wsAll.Range("F16").Value = 84
$wsAll.Range("F18").Value = 224
$wsAll.Range("F19").Value = 199
$wsAll.Range("F20").Value = 60
$wsAll.Range("F21").Value = 227
$wsAll.Range("F22").Value = 965
